# "Test case senarios and cases.xlsx" — widen column B and move the
# selection to C3 (single cell), matching the resize + re-click done in
# Excel before the file was re-uploaded.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Widen column B (was ~21.33 "characters" of stored width, now ~50.22).
$ws.Columns.Item(2).ColumnWidth = 49.333333333333336

# Collapse the old A1:C3 selection down to the single cell C3.
$ws.Range("C3").Select()
